$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the statement that relates to Report writing with "Report" in column E
# (written first so it lands before "File" in the shared-strings table)
$ws.Range("E48").Value = "Report"

# Mark the statements that relate to File I/O with "File" in column E
$ws.Range("E7").Value = "File"
$ws.Range("E10").Value = "File"
$ws.Range("E38").Value = "File"
$ws.Range("E43").Value = "File"
$ws.Range("E44").Value = "File"

# Restore selection to match the author's last-edited cell
$ws.Range("E10").Select() | Out-Null
